$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the F-column (time_taken) timestamps on the "data" sheet ---
# The panel was re-queried, so every row gets a freshly generated timestamp
# (same gene rows/values, later query run).
$timestamps = @{
  2 = "2021-10-05 14:34:15.457752"
  3 = "2021-10-05 14:34:15.457760"
  4 = "2021-10-05 14:34:15.457763"
  5 = "2021-10-05 14:34:15.457765"
  6 = "2021-10-05 14:34:15.457768"
  7 = "2021-10-05 14:34:15.457771"
  8 = "2021-10-05 14:34:15.457773"
  9 = "2021-10-05 14:34:15.457776"
  10 = "2021-10-05 14:34:15.457778"
  11 = "2021-10-05 14:34:15.457781"
  12 = "2021-10-05 14:34:15.457783"
  13 = "2021-10-05 14:34:15.457786"
  14 = "2021-10-05 14:34:15.457788"
  15 = "2021-10-05 14:34:15.457791"
  16 = "2021-10-05 14:34:15.457793"
  17 = "2021-10-05 14:34:15.457796"
  18 = "2021-10-05 14:34:15.457799"
  19 = "2021-10-05 14:34:15.457801"
  20 = "2021-10-05 14:34:15.457804"
  21 = "2021-10-05 14:34:15.457806"
  22 = "2021-10-05 14:34:15.457809"
  23 = "2021-10-05 14:34:15.457811"
  24 = "2021-10-05 14:34:15.457814"
  25 = "2021-10-05 14:34:15.457817"
  26 = "2021-10-05 14:34:15.457819"
  27 = "2021-10-05 14:34:15.457822"
  28 = "2021-10-05 14:34:15.457824"
  29 = "2021-10-05 14:34:15.457827"
  30 = "2021-10-05 14:34:15.457829"
  31 = "2021-10-05 14:34:15.457831"
  32 = "2021-10-05 14:34:15.457834"
  33 = "2021-10-05 14:34:15.457836"
  34 = "2021-10-05 14:34:15.457839"
  35 = "2021-10-05 14:34:15.457842"
  36 = "2021-10-05 14:34:15.457844"
  37 = "2021-10-05 14:34:15.457847"
  38 = "2021-10-05 14:34:15.457849"
  39 = "2021-10-05 14:34:15.457851"
  40 = "2021-10-05 14:34:15.457854"
  41 = "2021-10-05 14:34:15.457856"
  42 = "2021-10-05 14:34:15.457859"
  43 = "2021-10-05 14:34:15.457862"
  44 = "2021-10-05 14:34:15.457864"
  45 = "2021-10-05 14:34:15.457866"
  46 = "2021-10-05 14:34:15.457869"
  47 = "2021-10-05 14:34:15.457871"
  48 = "2021-10-05 14:34:15.457874"
  49 = "2021-10-05 14:34:15.457876"
  50 = "2021-10-05 14:34:15.457878"
  51 = "2021-10-05 14:34:15.457881"
  52 = "2021-10-05 14:34:15.457883"
  53 = "2021-10-05 14:34:15.457886"
  54 = "2021-10-05 14:34:15.457888"
  55 = "2021-10-05 14:34:15.457891"
  56 = "2021-10-05 14:34:15.457894"
  57 = "2021-10-05 14:34:15.457896"
  58 = "2021-10-05 14:34:15.457899"
  59 = "2021-10-05 14:34:15.457901"
  60 = "2021-10-05 14:34:15.457903"
  61 = "2021-10-05 14:34:15.457906"
  62 = "2021-10-05 14:34:15.457908"
  63 = "2021-10-05 14:34:15.457911"
  64 = "2021-10-05 14:34:15.457913"
  65 = "2021-10-05 14:34:15.457916"
  66 = "2021-10-05 14:34:15.457919"
  67 = "2021-10-05 14:34:15.457922"
}
foreach ($row in $timestamps.Keys) {
  $dataSheet.Cells.Item([int]$row, 6).Value2 = $timestamps[$row]
}

# --- Add a new "metadata" tab after "data", describing the panel export run ---
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# Match the outline summary settings used on the "data" sheet.
$metaSheet.Outline.SummaryRow = 1
$metaSheet.Outline.SummaryColumn = 1

# Copy the header/index cell formatting (bold, bordered, centered) from "data"
# so the new sheet matches its visual style, then fill in the text.
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122) | Out-Null
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Header row
$metaSheet.Cells.Item(1, 2).Value2 = "data_name"
$metaSheet.Cells.Item(1, 3).Value2 = "data_id"
$metaSheet.Cells.Item(1, 4).Value2 = "data_version"
$metaSheet.Cells.Item(1, 5).Value2 = "data_version_created"
$metaSheet.Cells.Item(1, 6).Value2 = "panel_query_time"
$metaSheet.Cells.Item(1, 7).Value2 = "panel_get_request"

# Data row
$metaSheet.Cells.Item(2, 1).Value2 = 0
$metaSheet.Cells.Item(2, 2).Value2 = "Joubert syndrome and other neurological ciliopathies"
$metaSheet.Cells.Item(2, 3).Value2 = 129

# data_version ("1.14") must stay a text value, not be coerced to a number
$versionCell = $metaSheet.Cells.Item(2, 4)
$versionCell.NumberFormat = "@"
$versionCell.Value2 = "1.14"
$versionCell.Style = "Normal"

$metaSheet.Cells.Item(2, 5).Value2 = "2021-09-06T05:27:41.902858Z"
$metaSheet.Cells.Item(2, 6).Value2 = "2021-10-05 14:34:15.454191"
$metaSheet.Cells.Item(2, 7).Value2 = "https://panelapp.agha.umccr.org/api/v1/panels/129/?format=json"

# Leave "data" as the active/selected sheet (unchanged from the original workbook view).
$dataSheet.Activate()
$dataSheet.Range("A1").Select() | Out-Null
